# Generate Report for Handoff
# Updates the localization-status workbook with the latest handoff
# information for the five files that just got handed off
# (600a1aba, 8ced41c0, 918f6229, b49c5156, cba2811d, d844475a -> rows 5-10).

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
# Column G = "Latest HO Xliff Generate Date" for rows 5-10
$overview = $wb.Worksheets.Item("Overview")
for ($row = 5; $row -le 10; $row++) {
    $overview.Range("G$row").Value2 = "2016-08-31 13:20:34"
}

# --- zh-cn sheet ------------------------------------------------------
# Column E = "Priority", Column H = "Latest Handoff Datetime" for rows 5-10
$zhcn = $wb.Worksheets.Item("zh-cn")
for ($row = 5; $row -le 10; $row++) {
    $zhcn.Range("E$row").Value2 = "ht"
    $zhcn.Range("H$row").Value2 = "2016-08-31 13:20:13"
}

# --- de-de sheet -------------------------------------------------------
# Column E = "Priority", Column H = "Latest Handoff Datetime" for rows 5-10
$dede = $wb.Worksheets.Item("de-de")
for ($row = 5; $row -le 10; $row++) {
    $dede.Range("E$row").Value2 = "ht"
    $dede.Range("H$row").Value2 = "2016-08-31 13:20:34"
}
